$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.368.61"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.884.54"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7135"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.34"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08098"
$ws.Range("E8").Value = "  +4.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3137"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.34"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08359"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").Value = "1.875.94"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7223"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.252"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.15"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.288"
$ws.Range("E16").Value = "  +4.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008473"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").Value = "29.375.76"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.66"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").Value = "2.119.92"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.813"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.33"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.430"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.359"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.219"
$ws.Range("E32").Value = "  -3.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05386"
$ws.Range("E33").Value = "  +2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.957"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7529"
$ws.Range("E35").Value = "  +0.94%  "

$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01883"
$ws.Range("E38").Value = "  +1.19%  "

$ws.Range("D39").Value = "1.285.44"
$ws.Range("E39").Value = "  +10.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.745"
$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.571"
$ws.Range("E41").Value = "  +3.15%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.62"
$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8933"
$ws.Range("E43").Value = "  +0.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "110.39"
$ws.Range("E44").Value = "  +3.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000129"
$ws.Range("E46").Value = "  +6.89%  "

$ws.Range("D47").Value = "2.020.81"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.806"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5214"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.494"
$ws.Range("E50").Value = "  +1.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4368"
$ws.Range("E51").Value = "  +1.49%  "
